$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-11 15:14:16", 0.0008),
    @("2023-12-11 15:14:34", 0.0014),
    @("2023-12-11 15:14:46", 0.0004),
    @("2023-12-11 15:15:08", 0.0012)
)

$startRow = 191
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
